# This script applies scheduled market-price/profit updates to the
# Leve profit tables (ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR) produced by the
# automated runner, matching the refreshed currentAveragePrice* /
# LevePrice* / LeveProfit* figures (columns H-N) for the affected rows.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 19
$ws.Range("H19").Value = 983.1875
$ws.Range("I19").Value = 1121.5714
$ws.Range("J19").Value = 875.55554
$ws.Range("K19").Value = 1121.5714
$ws.Range("L19").Value = 875.55554
$ws.Range("M19").Value = -946.5714
$ws.Range("N19").Value = -1225.55554
# Row 33
$ws.Range("H33").Value = 83.77778000000001
$ws.Range("I33").Value = 83.77778000000001
$ws.Range("K33").Value = 83.77778000000001
$ws.Range("M33").Value = 145.22222
# Row 55
$ws.Range("H55").Value = 99.083336
$ws.Range("I55").Value = 98.77778000000001
$ws.Range("J55").Value = 100
$ws.Range("K55").Value = 98.77778000000001
$ws.Range("L55").Value = 100
$ws.Range("M55").Value = 115.22222
$ws.Range("N55").Value = -528
# Row 70
$ws.Range("H70").Value = 2628.3713
$ws.Range("I70").Value = 3666
$ws.Range("J70").Value = 2213.32
$ws.Range("K70").Value = 10998
$ws.Range("L70").Value = 6639.960000000001
$ws.Range("M70").Value = -10728
$ws.Range("N70").Value = -7179.960000000001
# Row 73
$ws.Range("H73").Value = 2628.3713
$ws.Range("I73").Value = 3666
$ws.Range("J73").Value = 2213.32
$ws.Range("K73").Value = 10998
$ws.Range("L73").Value = 6639.960000000001
$ws.Range("M73").Value = -10062
$ws.Range("N73").Value = -8511.960000000001
# Row 86
$ws.Range("H86").Value = 2510.3333
$ws.Range("I86").Value = 1048.25
$ws.Range("J86").Value = 3680
$ws.Range("K86").Value = 1048.25
$ws.Range("L86").Value = 3680
$ws.Range("M86").Value = 74.75
$ws.Range("N86").Value = -5926
# Row 87
$ws.Range("H87").Value = 39106.75
$ws.Range("J87").Value = 39106.75
$ws.Range("L87").Value = 39106.75
$ws.Range("N87").Value = -41602.75
# Row 89
$ws.Range("H89").Value = 2510.3333
$ws.Range("I89").Value = 1048.25
$ws.Range("J89").Value = 3680
$ws.Range("K89").Value = 5241.25
$ws.Range("L89").Value = 18400
$ws.Range("M89").Value = 374.75
$ws.Range("N89").Value = -29632
# Row 90
$ws.Range("H90").Value = 39106.75
$ws.Range("J90").Value = 39106.75
$ws.Range("L90").Value = 117320.25
$ws.Range("N90").Value = -129800.25

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 36
$ws.Range("H36").Value = 0
$ws.Range("I36").Value = 0
$ws.Range("K36").Value = 0
$ws.Range("M36").ClearContents()
# Row 40
$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").ClearContents()
# Row 44
$ws.Range("H44").Value = 4800
$ws.Range("I44").Value = 4800
$ws.Range("J44").Value = 0
$ws.Range("K44").Value = 4800
$ws.Range("L44").Value = 0
$ws.Range("M44").Value = -4358
$ws.Range("N44").ClearContents()
# Row 48
$ws.Range("H48").Value = 5909.2
$ws.Range("I48").Value = 46
$ws.Range("J48").Value = 7375
$ws.Range("K48").Value = 46
$ws.Range("L48").Value = 7375
$ws.Range("M48").Value = 430
$ws.Range("N48").Value = -8327
# Row 50
$ws.Range("H50").Value = 22445
$ws.Range("I50").Value = 19990
$ws.Range("J50").Value = 24900
$ws.Range("K50").Value = 19990
$ws.Range("L50").Value = 24900
$ws.Range("M50").Value = -19365
$ws.Range("N50").Value = -26150
# Row 51
$ws.Range("H51").Value = 27595
$ws.Range("I51").Value = 24400
$ws.Range("J51").Value = 28660
$ws.Range("K51").Value = 24400
$ws.Range("L51").Value = 28660
$ws.Range("M51").Value = -23664
$ws.Range("N51").Value = -30132
# Row 56
$ws.Range("H56").Value = 10500
$ws.Range("I56").Value = 10500
$ws.Range("J56").Value = 0
$ws.Range("K56").Value = 10500
$ws.Range("L56").Value = 0
$ws.Range("M56").Value = -9655
$ws.Range("N56").ClearContents()
# Row 57
$ws.Range("H57").Value = 9000
$ws.Range("J57").Value = 9000
$ws.Range("L57").Value = 9000
$ws.Range("N57").Value = -10120
# Row 61
$ws.Range("H61").Value = 27595
$ws.Range("I61").Value = 24400
$ws.Range("J61").Value = 28660
$ws.Range("K61").Value = 24400
$ws.Range("L61").Value = 28660
$ws.Range("M61").Value = -24052
$ws.Range("N61").Value = -29356
# Row 62
$ws.Range("H62").Value = 4849.1665
$ws.Range("I62").Value = 5368.421
$ws.Range("J62").Value = 2876
$ws.Range("K62").Value = 5368.421
$ws.Range("L62").Value = 2876
$ws.Range("M62").Value = -4744.421
$ws.Range("N62").Value = -4124
# Row 65
$ws.Range("H65").Value = 4849.1665
$ws.Range("I65").Value = 5368.421
$ws.Range("J65").Value = 2876
$ws.Range("K65").Value = 26842.105
$ws.Range("L65").Value = 14380
$ws.Range("M65").Value = -23722.105
$ws.Range("N65").Value = -20620
# Row 134
$ws.Range("H134").Value = 3354.3157
$ws.Range("I134").Value = 5043.3335
$ws.Range("J134").Value = 1834.2
$ws.Range("K134").Value = 15130.0005
$ws.Range("L134").Value = 5502.6
$ws.Range("M134").Value = -12595.0005
$ws.Range("N134").Value = -10572.6

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 47
$ws.Range("H47").Value = 287.375
$ws.Range("I47").Value = 287.375
$ws.Range("J47").Value = 0
$ws.Range("K47").Value = 862.125
$ws.Range("L47").Value = 0
$ws.Range("M47").Value = -431.125
$ws.Range("N47").ClearContents()
# Row 68
$ws.Range("H68").Value = 1946.75
$ws.Range("I68").Value = 860
$ws.Range("J68").Value = 2598.8
$ws.Range("K68").Value = 2580
$ws.Range("L68").Value = 7796.400000000001
$ws.Range("M68").Value = -1769
$ws.Range("N68").Value = -9418.400000000001
# Row 71
$ws.Range("H71").Value = 1946.75
$ws.Range("I71").Value = 860
$ws.Range("J71").Value = 2598.8
$ws.Range("K71").Value = 7740
$ws.Range("L71").Value = 23389.2
$ws.Range("M71").Value = -3684
$ws.Range("N71").Value = -31501.2
# Row 92
$ws.Range("H92").Value = 0
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("M92").ClearContents()
$ws.Range("N92").ClearContents()
# Row 131
$ws.Range("H131").Value = 821.1395
$ws.Range("I131").Value = 235.44444
$ws.Range("J131").Value = 976.17645
$ws.Range("K131").Value = 706.33332
$ws.Range("L131").Value = 2928.52935
$ws.Range("M131").Value = 4333.66668
$ws.Range("N131").Value = -13008.52935

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 46
$ws.Range("H46").Value = 6416.6665
$ws.Range("I46").Value = 1625
$ws.Range("J46").Value = 16000
$ws.Range("K46").Value = 1625
$ws.Range("L46").Value = 16000
$ws.Range("M46").Value = -1469
$ws.Range("N46").Value = -16312
# Row 58
$ws.Range("H58").Value = 13950
$ws.Range("J58").Value = 13950
$ws.Range("L58").Value = 13950
$ws.Range("N58").Value = -14504
# Row 80
$ws.Range("H80").Value = 3019.125
$ws.Range("I80").Value = 2163
$ws.Range("J80").Value = 4446
$ws.Range("K80").Value = 2163
$ws.Range("L80").Value = 4446
$ws.Range("M80").Value = -1165
$ws.Range("N80").Value = -6442
# Row 83
$ws.Range("H83").Value = 3019.125
$ws.Range("I83").Value = 2163
$ws.Range("J83").Value = 4446
$ws.Range("K83").Value = 10815
$ws.Range("L83").Value = 22230
$ws.Range("M83").Value = -5823
$ws.Range("N83").Value = -32214
# Row 132
$ws.Range("H132").Value = 37937.71
$ws.Range("I132").Value = 42900.074
$ws.Range("J132").Value = 4441.75
$ws.Range("K132").Value = 128700.222
$ws.Range("L132").Value = 13325.25
$ws.Range("M132").Value = -126170.222
$ws.Range("N132").Value = -18385.25

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 37
$ws.Range("H37").Value = 5000
$ws.Range("I37").Value = 5000
$ws.Range("J37").Value = 0
$ws.Range("K37").Value = 5000
$ws.Range("L37").Value = 0
$ws.Range("M37").Value = -4797
$ws.Range("N37").ClearContents()
# Row 40
$ws.Range("H40").Value = 6500
$ws.Range("J40").Value = 6500
$ws.Range("L40").Value = 6500
$ws.Range("N40").Value = -6798
# Row 52
$ws.Range("H52").Value = 8173.5
$ws.Range("I52").Value = 3000
$ws.Range("J52").Value = 13347
$ws.Range("K52").Value = 3000
$ws.Range("L52").Value = 13347
$ws.Range("M52").Value = -2774
$ws.Range("N52").Value = -13799
# Row 58
$ws.Range("H58").Value = 17088
$ws.Range("I58").Value = 15585
$ws.Range("J58").Value = 20094
$ws.Range("K58").Value = 15585
$ws.Range("L58").Value = 20094
$ws.Range("M58").Value = -15277
$ws.Range("N58").Value = -20710

